$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 128
$ws.Range("B2").Value = 0.39
$ws.Range("B3").Value = 86

$ws.Range("B3").Select()
